$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9, shifting existing rows 9.. down by one.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record.
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44749
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100108
$ws.Cells.Item(9, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(9, 9).Value = 100108007
$ws.Cells.Item(9, 10).Value = "Coco"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 14).Value = 28000
$ws.Cells.Item(9, 15).Value = 30000
$ws.Cells.Item(9, 16).Value = 28800
$ws.Cells.Item(9, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(9, 18).Value = "Perú"
$ws.Cells.Item(9, 19).Value = 1440
$ws.Cells.Item(9, 20).Value = 20
